# Adds a new "2022-Q3" quarterly sheet (right after "总计") with the latest
# fund-holding snapshot, and records the new quarter in the "总计" rollup
# sheet. All the other quarterly tabs (2022-Q2 .. 2020-Q4) are left exactly
# as they are -- they simply shift one tab to the right, which Excel does
# for us automatically once the new sheet is inserted ahead of them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the "2022-Q3" sheet by duplicating "2022-Q2" (so it inherits
#    the same layout/formatting) and placing the copy right before it.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)

$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
}

# Existing three funds -- refreshed figures for 2022-Q3.
Set-TextValue $q3.Range("C2") "华宝标普美国品质消费股票（LOF）美元"
Set-TextValue $q3.Range("D2") "3.59"
Set-TextValue $q3.Range("E2") "94.37"
Set-TextValue $q3.Range("F2") "3.10"
Set-TextValue $q3.Range("G2") "0.1113"

Set-TextValue $q3.Range("C3") "华宝标普美国品质消费股票（LOF）人民币A"
Set-TextValue $q3.Range("D3") "2.86"
Set-TextValue $q3.Range("E3") "94.37"
Set-TextValue $q3.Range("F3") "3.10"
Set-TextValue $q3.Range("G3") "0.0887"

Set-TextValue $q3.Range("C4") "华宝标普美国品质消费股票（LOF）人民币C"
Set-TextValue $q3.Range("D4") "0.73"
Set-TextValue $q3.Range("E4") "94.37"
Set-TextValue $q3.Range("F4") "3.10"
Set-TextValue $q3.Range("G4") "0.0226"

# Two brand-new funds showed up in the 2022-Q3 snapshot -- add rows 5 & 6,
# copying row 4's formatting first so borders/styles stay consistent.
$q3.Range("A4:H4").Copy()
$q3.Range("A5:H5").PasteSpecial(-4122)
$q3.Range("A4:H4").Copy()
$q3.Range("A6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q3.Range("A5").Value2 = 3
Set-TextValue $q3.Range("B5") "011706"
Set-TextValue $q3.Range("C5") "长信美国标准普尔100等权重指数增强（QDII）美元"
Set-TextValue $q3.Range("D5") "0.39"
Set-TextValue $q3.Range("E5") "82.64"
Set-TextValue $q3.Range("F5") "0.85"
Set-TextValue $q3.Range("G5") "0.0033"
$q3.Range("H5").Value2 = 9

$q3.Range("A6").Value2 = 4
Set-TextValue $q3.Range("B6") "519981"
Set-TextValue $q3.Range("C6") "长信美国标准普尔100等权重指数增强（QDII）人民币"
Set-TextValue $q3.Range("D6") "0.39"
Set-TextValue $q3.Range("E6") "82.64"
Set-TextValue $q3.Range("F6") "0.85"
Set-TextValue $q3.Range("G6") "0.0033"
$q3.Range("H6").Value2 = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" rollup sheet: a new row for 2022-Q3 goes in right
#    under the header, every other quarter shifts down one row, and the
#    row-index column (A) is renumbered 0..7.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Give the brand-new last row (for 2020-Q4) the same "A column" style the
# rest of the table already uses.
$zj.Range("A8").Copy()
$zj.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @(0, "2022-Q3", 5, 0.23),
    @(1, "2022-Q2", 3, 0.17),
    @(2, "2022-Q1", 3, 0.19),
    @(3, "2021-Q4", 3, 0.26),
    @(4, "2021-Q3", 3, 0.23),
    @(5, "2021-Q2", 3, 0.21),
    @(6, "2021-Q1", 3, 0.2),
    @(7, "2020-Q4", 3, 0.18)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $zj.Cells.Item($r, 1).Value2 = $row[0]
    $zj.Cells.Item($r, 2).Value2 = $row[1]
    $zj.Cells.Item($r, 3).Value2 = $row[2]
    $zj.Cells.Item($r, 4).Value2 = $row[3]
}

Write-Output "2022-Q3 used range: $($q3.UsedRange.Address())"
Write-Output "总计 used range: $($zj.UsedRange.Address())"
